$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.829.54'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.634.64'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = "'215.19"
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").Value = "'0.5033"
$ws.Range("E6").Value = '  -2.32%  '
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").Value = "'0.2575"
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").Value = "'0.06412"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").Value = "'0.07674"
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("D12").Value = '1.640.29'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").Value = "'4.240"
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("D14").Value = '1.858.20'
$ws.Range("E14").Value = '  -1.56%  '
$ws.Range("D15").Value = "'0.5461"
$ws.Range("D16").Value = '0.0₅7928'
$ws.Range("E16").Value = '  -1.84%  '
$ws.Range("D17").Value = "'63.55"
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").Value = '25.844.58'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").Value = "'203.10"
$ws.Range("E20").Value = '  -4.09%  '
$ws.Range("D21").Value = "'4.317"
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("D22").Value = "'9.937"
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = "'5.967"
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").Value = "'1.936"
$ws.Range("E25").Value = '  +10.29%  '
$ws.Range("D26").Value = "'140.88"
$ws.Range("E26").Value = '  -2.54%  '
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("D28").Value = "'15.72"
$ws.Range("D29").Value = "'6.695"
$ws.Range("E29").Value = '  -4.22%  '
$ws.Range("D30").Value = "'1.240"
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("D31").Value = "'0.04974"
$ws.Range("E31").Value = '  -4.70%  '
$ws.Range("D32").Value = "'3.276"
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("D33").Value = "'3.184"
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("D34").Value = "'1.535"
$ws.Range("E34").Value = '  -2.48%  '
$ws.Range("D35").Value = "'2.351"
$ws.Range("E35").Value = '  -0.82%  '
$ws.Range("D36").Value = '1.176.01'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").Value = "'0.8925"
$ws.Range("E37").Value = '  -4.38%  '
$ws.Range("D38").Value = "'2.617"
$ws.Range("E38").Value = '  -5.17%  '
$ws.Range("D39").Value = "'0.5583"
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("D40").Value = "'0.01557"
$ws.Range("E40").Value = '  -2.37%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("B42").Value = 'mCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D42").Value = "'2.543"
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("D43").Value = "'5.650"
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = "'0.8042"
$ws.Range("E44").Value = '  -4.98%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = "'99.31"
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").Value = '1.770.27'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").Value = "'0.4511"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = "'1.006"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").Value = "'54.82"
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").Value = "'0.05031"
$ws.Range("E51").Value = '  -0.67%  '
